# "error solve ifrs list" - replace the (incorrectly scaled) financial
# figures in rows 2-6 with the corrected values, and wipe out the stale
# data rows 7-9 (their annual columns, D:AJ, are cleared entirely while
# the row stays with its A/B/C labels).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (연간 2014/12)
$ws.Range("D2").Value = 1122
$ws.Range("E2").Value = -10
$ws.Range("F2").Value = -10
$ws.Range("G2").Value = 5
$ws.Range("H2").Value = 2
$ws.Range("I2").Value = 4
$ws.Range("J2").Value = -2
$ws.Range("K2").Value = 1568
$ws.Range("L2").Value = 212
$ws.Range("M2").Value = 1356
$ws.Range("N2").Value = 1312
$ws.Range("O2").Value = 44
$ws.Range("P2").Value = 60
$ws.Range("Q2").Value = 108
$ws.Range("R2").Value = -109
$ws.Range("S2").Value = -28
$ws.Range("T2").Value = 5
$ws.Range("U2").Value = 104
$ws.Range("V2").Value = 51
$ws.Range("W2").Value = -0.85
$ws.Range("X2").Value = 0.17
$ws.Range("Y2").Value = 0.3
$ws.Range("Z2").Value = 0.12
$ws.Range("AA2").Value = 15.66
$ws.Range("AB2").Value = 2195.64
$ws.Range("AC2").Value = 34
$ws.Range("AD2").Value = 475.12
$ws.Range("AE2").Value = 13234
$ws.Range("AF2").Value = 1.23
$ws.Range("AG2").Value = 200
$ws.Range("AH2").Value = 1.23
$ws.Range("AI2").Value = 493.92
$ws.Range("AJ2").Value = 11700000

# Row 3 (연간 2015/12)
$ws.Range("D3").Value = 1062
$ws.Range("E3").Value = 62
$ws.Range("F3").Value = 62
$ws.Range("G3").Value = 54
$ws.Range("H3").Value = 45
$ws.Range("I3").Value = 39
$ws.Range("J3").Value = 5
$ws.Range("K3").Value = 1632
$ws.Range("L3").Value = 246
$ws.Range("M3").Value = 1386
$ws.Range("N3").Value = 1337
$ws.Range("O3").Value = 49
$ws.Range("P3").Value = 60
$ws.Range("Q3").Value = 181
$ws.Range("R3").Value = -133
$ws.Range("S3").Value = -30
$ws.Range("T3").Value = 5
$ws.Range("U3").Value = 176
$ws.Range("V3").Value = 40
$ws.Range("W3").Value = 5.89
$ws.Range("X3").Value = 4.19
$ws.Range("Y3").Value = 2.98
$ws.Range("Z3").Value = 2.78
$ws.Range("AA3").Value = 17.79
$ws.Range("AB3").Value = 2225.33
$ws.Range("AC3").Value = 337
$ws.Range("AD3").Value = 30.69
$ws.Range("AE3").Value = 13490
$ws.Range("AF3").Value = 0.77
$ws.Range("AG3").Value = 200
$ws.Range("AH3").Value = 1.93
$ws.Range("AI3").Value = 50.25
$ws.Range("AJ3").Value = 11700000

# Row 4 (연간 2016/12)
$ws.Range("D4").Value = 1285
$ws.Range("E4").Value = 194
$ws.Range("F4").Value = 194
$ws.Range("G4").Value = 269
$ws.Range("H4").Value = 208
$ws.Range("I4").Value = 196
$ws.Range("J4").Value = 13
$ws.Range("K4").Value = 1820
$ws.Range("L4").Value = 275
$ws.Range("M4").Value = 1545
$ws.Range("N4").Value = 1483
$ws.Range("O4").Value = 62
$ws.Range("P4").Value = 60
$ws.Range("Q4").Value = 159
$ws.Range("R4").Value = -86
$ws.Range("S4").Value = -66
$ws.Range("T4").Value = 20
$ws.Range("U4").Value = 140
$ws.Range("V4").Value = 31
$ws.Range("W4").Value = 15.12
$ws.Range("X4").Value = 16.21
$ws.Range("Y4").Value = 13.87
$ws.Range("Z4").Value = 12.07
$ws.Range("AA4").Value = 17.8
$ws.Range("AB4").Value = 2528.88
$ws.Range("AC4").Value = 1672
$ws.Range("AD4").Value = 7.57
$ws.Range("AE4").Value = 15512
$ws.Range("AF4").Value = 0.82
$ws.Range("AG4").Value = 400
$ws.Range("AH4").Value = 3.16
$ws.Range("AI4").Value = 19.55
$ws.Range("AJ4").Value = 11700000

# Row 5 (연간 2017/12)
$ws.Range("D5").Value = 1336
$ws.Range("E5").Value = 166
$ws.Range("F5").Value = 166
$ws.Range("G5").Value = 177
$ws.Range("H5").Value = 141
$ws.Range("I5").Value = 130
$ws.Range("J5").Value = 10
$ws.Range("K5").Value = 1862
$ws.Range("L5").Value = 251
$ws.Range("M5").Value = 1611
$ws.Range("N5").Value = 1546
$ws.Range("O5").Value = 66
$ws.Range("P5").Value = 60
$ws.Range("Q5").Value = 102
$ws.Range("R5").Value = -62
$ws.Range("S5").Value = -74
$ws.Range("T5").Value = 33
$ws.Range("U5").Value = 69
$ws.Range("V5").Value = 30
$ws.Range("W5").Value = 12.39
$ws.Range("X5").Value = 10.52
$ws.Range("Y5").Value = 8.6
$ws.Range("Z5").Value = 7.64
$ws.Range("AA5").Value = 15.55
$ws.Range("AB5").Value = 2682.69
$ws.Range("AC5").Value = 1113
$ws.Range("AD5").Value = 11.1
$ws.Range("AE5").Value = 16565
$ws.Range("AF5").Value = 0.75
$ws.Range("AG5").Value = 460
$ws.Range("AH5").Value = 3.72
$ws.Range("AI5").Value = 32.96
$ws.Range("AJ5").Value = 11700000

# Row 6 (연간 2018/12) - J6/O6 were already blank before and stay blank
$ws.Range("D6").Value = 1221
$ws.Range("E6").Value = 85
$ws.Range("F6").Value = 85
$ws.Range("G6").Value = 89
$ws.Range("H6").Value = 72
$ws.Range("I6").Value = 67
$ws.Range("J6").ClearContents()
$ws.Range("K6").Value = 1879
$ws.Range("L6").Value = 250
$ws.Range("M6").Value = 1628
$ws.Range("N6").Value = 1562
$ws.Range("O6").ClearContents()
$ws.Range("P6").Value = 60
$ws.Range("Q6").Value = 133
$ws.Range("R6").Value = -76
$ws.Range("S6").Value = -40
$ws.Range("T6").Value = 38
$ws.Range("U6").Value = 96
$ws.Range("V6").Value = 37
$ws.Range("W6").Value = 6.94
$ws.Range("X6").Value = 5.89
$ws.Range("Y6").Value = 4.29
$ws.Range("Z6").Value = 3.84
$ws.Range("AA6").Value = 15.37
$ws.Range("AB6").Value = 2712.73
$ws.Range("AC6").Value = 570
$ws.Range("AD6").Value = 19.12
$ws.Range("AE6").Value = 16738
$ws.Range("AF6").Value = 0.65
$ws.Range("AG6").Value = 250
$ws.Range("AH6").Value = 2.29
$ws.Range("AI6").Value = 34.98
$ws.Range("AJ6").Value = 11700000

# Rows 7-9 (연간 2019/12(E), 2020/12(E), 2021/12(E)) - these forecast
# rows no longer have data; clear D:AJ but keep the A/B/C row labels.
$ws.Range("D7:AJ7").ClearContents()
$ws.Range("D8:AJ8").ClearContents()
$ws.Range("D9:AJ9").ClearContents()
